$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.949.72'
$ws.Range('E2').Value = '  +0.62%  '
$ws.Range('D3').Value = '2.298.13'
$ws.Range('E3').Value = '  +0.18%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').Value = '300.10'
$ws.Range('E5').Value = '  -0.43%  '
$ws.Range('D6').Value = '97.05'
$ws.Range('E6').Value = '  -1.06%  '
$ws.Range('D7').Value = '0.507'
$ws.Range('E7').Value = '  +0.26%  '
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('E9').Value = '  -0.58%  '
$ws.Range('D10').Value = '33.56'
$ws.Range('E10').Value = '  -2.95%  '
$ws.Range('D11').Value = '0.0793'
$ws.Range('E11').Value = '  +0.21%  '
$ws.Range('D12').Value = '49.21'
$ws.Range('E12').Value = '  -3.88%  '
$ws.Range('E13').Value = '  +2.42%  '
$ws.Range('D14').Value = '17.02'
$ws.Range('E14').Value = '  +10.38%  '
$ws.Range('D15').Value = '6.77'
$ws.Range('E15').Value = '  +0.58%  '
$ws.Range('D16').Value = '2.659.81'
$ws.Range('E16').Value = '  +0.01%  '
$ws.Range('D17').Value = '2.288.66'
$ws.Range('E17').Value = '  +0.00%  '
$ws.Range('D18').Value = '0.809'
$ws.Range('E18').Value = '  +2.11%  '
$ws.Range('D19').Value = '42.939.21'
$ws.Range('E19').Value = '  +0.56%  '
$ws.Range('D20').Value = '0.0₃0901'
$ws.Range('E20').Value = '  +0.39%  '
$ws.Range('D21').Value = '11.57'
$ws.Range('E21').Value = '  -0.40%  '
$ws.Range('E22').Value = '  +0.50%  '
$ws.Range('D23').Value = '67.83'
$ws.Range('E23').Value = '  +0.98%  '
$ws.Range('D24').Value = '236.35'
$ws.Range('E24').Value = '  +0.60%  '
$ws.Range('D25').Value = '2.03'
$ws.Range('E25').Value = '  +4.02%  '
$ws.Range('E26').Value = '  +0.01%  '
$ws.Range('D27').Value = '2.45'
$ws.Range('E27').Value = '  -1.97%  '
$ws.Range('D28').Value = '24.45'
$ws.Range('E28').Value = '  -0.69%  '
$ws.Range('D29').Value = '166.34'
$ws.Range('E29').Value = '  +1.04%  '
$ws.Range('D30').Value = '33.94'
$ws.Range('E30').Value = '  -0.79%  '
$ws.Range('E31').Value = '  -5.65%  '
$ws.Range('D32').Value = '9.12'
$ws.Range('E32').Value = '  +0.04%  '
$ws.Range('E33').Value = '  -0.09%  '
$ws.Range('D34').Value = '4.68'
$ws.Range('E34').Value = '  +6.97%  '
$ws.Range('D35').Value = '4.94'
$ws.Range('E35').Value = '  -0.84%  '
$ws.Range('E36').Value = '  -1.00%  '
$ws.Range('D37').Value = '16.82'
$ws.Range('E37').Value = '  +3.60%  '
$ws.Range('E38').Value = '  -0.75%  '
$ws.Range('E39').Value = '  -0.30%  '
$ws.Range('E40').Value = '  +0.46%  '
$ws.Range('D41').Value = '1.76'
$ws.Range('E41').Value = '  -1.41%  '
$ws.Range('D42').Value = '0.109'
$ws.Range('E42').Value = '  -0.74%  '
$ws.Range('D43').Value = '2.41'
$ws.Range('E43').Value = '  +0.17%  '
$ws.Range('D44').Value = '1.980.82'
$ws.Range('E44').Value = '  +0.40%  '
$ws.Range('E45').Value = '  -0.59%  '
$ws.Range('B46').Value = 'FraxShare'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D46').Value = '9.88'
$ws.Range('E46').Value = '  +0.67%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = '17.68'
$ws.Range('E47').Value = '  -1.73%  '
$ws.Range('D48').Value = '2.85'
$ws.Range('E48').Value = '  -0.58%  '
$ws.Range('D49').Value = '2.526.26'
$ws.Range('E49').Value = '  +0.04%  '
$ws.Range('D50').Value = '53.33'
$ws.Range('E50').Value = '  -0.30%  '
$ws.Range('D51').Value = '4.58'
$ws.Range('E51').Value = '  -3.52%  '
